# Testplan.docx - "Test 4: Aansturen van de waterpompen met Mosfets"
#
# 1. Replace the placeholder text under "Acceptatiecriteria" with the
#    actual acceptance-criteria text (with spell-check markers around the
#    loanwords "arduino" and "Mosfets", matching the author's edit).
# 2. Replace the placeholder text under "Waarnemingen" with a 2x2
#    observations table (header row "Test nr." / "Waarnemingen", data
#    row "1" / "-"), mirroring the table used for the other tests in
#    this document.

$d = $word.ActiveDocument

# Locate the "Test 4" section: from its heading up to (but not including)
# the next test's heading ("Test 5"). Using Find keeps this robust against
# any paragraph-count differences rather than hard-coding indices.
$startRange = $d.Content.Duplicate
$startRange.Find.Execute("Aansturen van de waterpompen met Mosfets", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$endRange = $d.Range($startRange.End, $d.Content.End)
$endRange.Find.Execute("Algoritme gebaseerd op sensordata", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$sectionEnd = $endRange.Start
$sectionStart = $startRange.End

# --- 1. "Acceptatiecriteria" placeholder -> real answer text ---------------
$critRange = $d.Range($sectionStart, $sectionEnd)
$critRange.Find.Execute("(Wanneer is de test voldaan*stellen.)", $true, $false, $true, $false, $false, $true, 1, $false, "", 0) | Out-Null

$critPara = $critRange.Paragraphs(1).Range
$critXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:r><w:t xml:space="preserve">De waterpompen worden correct aangestuurd op basis van de signalen die de </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>arduino</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> uitstuurt, door middel van de </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Mosfets</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">. </w:t></w:r>' +
    '</w:p>'
$critPara.InsertXML($critXml) | Out-Null

# --- 2. "Waarnemingen" placeholder -> observations table --------------------
# Recompute the section end, since the previous InsertXML may have shifted it.
$endRange2 = $d.Range($startRange.End, $d.Content.End)
$endRange2.Find.Execute("Algoritme gebaseerd op sensordata", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$obsRange = $d.Range($startRange.End, $endRange2.Start)
$obsRange.Find.Execute("(Feitelijke, objectieve observaties*trekken.)", $true, $false, $true, $false, $false, $true, 1, $false, "", 0) | Out-Null

$obsPara = $obsRange.Paragraphs(1).Range
$tblXml = '<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/>' +
    '<w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr>' +
    '<w:tblGrid><w:gridCol w:w="1345"/><w:gridCol w:w="7717"/></w:tblGrid>' +
    '<w:tr>' +
        '<w:tc><w:tcPr><w:tcW w:w="1345" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Test nr.</w:t></w:r></w:p></w:tc>' +
        '<w:tc><w:tcPr><w:tcW w:w="7717" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Waarnemingen</w:t></w:r></w:p></w:tc>' +
    '</w:tr>' +
    '<w:tr>' +
        '<w:tc><w:tcPr><w:tcW w:w="1345" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>1</w:t></w:r></w:p></w:tc>' +
        '<w:tc><w:tcPr><w:tcW w:w="7717" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>-</w:t></w:r></w:p></w:tc>' +
    '</w:tr>' +
    '</w:tbl>'
$obsPara.InsertXML($tblXml) | Out-Null

Write-Host "Test 4 Acceptatiecriteria/Waarnemingen updated."
